$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "U2"
$ws.Range("B24").Value = "E73-2G4M08S1C"
$ws.Range("C24").Value = "E73-2G4M08S1C-52840"
$ws.Range("D24").Value = 48260000
$ws.Range("E24").Value = -56515000
$ws.Range("F24").Value = 270000000
$ws.Range("G24").Value = "top"

$ws.Range("C25").Select()
